$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 171, pushing the existing rows 171-240 down to 172-241.
$ws.Rows("171:171").Insert()

# Populate the newly inserted row 171 with the new weekly record.
# Non-numeric/categorical columns mirror the rest of this data block.
$ws.Cells.Item(171, 1).Value = 8
$ws.Cells.Item(171, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(171, 3).Value = "Coquimbo"
$ws.Cells.Item(171, 4).Value = 44825
$ws.Cells.Item(171, 5).Value = 4
$ws.Cells.Item(171, 6).Value = 100112037
$ws.Cells.Item(171, 7).Value = "Cebollín"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 1300
$ws.Cells.Item(171, 11).Value = 1400
$ws.Cells.Item(171, 12).Value = 1600
$ws.Cells.Item(171, 13).Value = 1500
$ws.Cells.Item(171, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(171, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(171, 16).Value = 250
$ws.Cells.Item(171, 17).Value = 6
$ws.Cells.Item(171, 18).Value = "Hortaliza"
